# Cambiato modo di calcolare i valori medi per ecosistemi (approvvigionamento).
# Ora la media pesata avviene sui valori di eco_contribution per classe corine.
# I pesi sono le superfici corine.
# Add a new worksheet "Foglio1" after "Sheet1" with the new calculation.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Foglio1"

# Row 1: valori (unit values) for the 4 corine classes
$ws2.Range("A1").Value = 30
$ws2.Range("B1").Value = 45
$ws2.Range("C1").Value = 60
$ws2.Range("D1").Value = 70
$ws2.Range("E1").Value = "valori"

# Row 2: pesi superfici (corine surface weights)
$ws2.Range("A2").Value = 20
$ws2.Range("B2").Value = 10
$ws2.Range("C2").Value = 59
$ws2.Range("D2").Value = 63
$ws2.Range("E2").Value = "pesi superfici"

# Row 3: eco con (eco_contribution per unit)
$ws2.Range("A3").Value = 0.2
$ws2.Range("B3").Value = 0.36
$ws2.Range("C3").Value = 0.1
$ws2.Range("D3").Value = 0.15
$ws2.Range("E3").Value = "eco con"

# Row 5: controbution = valori * eco con
$ws2.Range("A5").Formula = "=+A1*A3"
$ws2.Range("B5:D5").Formula = "=+B1*B3"
$ws2.Range("E5").Value = "controbution"
$ws2.Range("G5").Formula = "=+AVERAGE(A5:D5)"

# Row 7: weighted average of unit values, weighted by corine surfaces
$ws2.Range("A7").Formula = "=+SUMPRODUCT(A1:D1,A2:D2)/SUM(A2:D2)"
$ws2.Range("B7").Value = "media ponderata valori unitari"

# Row 8: weighted average of eco con, weighted by corine surfaces (array formula)
$ws2.Range("A8").FormulaArray = "=+SUMPRODUCT(A3:D3,A2:D2)/SUM(A2:D2)"
$ws2.Range("B8").Value = "media ponderata eco con"

# Row 9: weighted average unit value * weighted average eco con
$ws2.Range("A9").Formula = "=+A7*A8"

$ws2.Range("G6").Select()
